$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 ("has_cv2" / duplicate CNS description) is removed entirely;
# all subsequent rows (14-27) shift up by one to become rows 13-26.
$ws.Rows("13:13").Delete()

# After the shift, update the renamed/rewritten rows near the bottom of the sheet.
# Former "FinalTx" (now row 24) becomes "FinalTx_coll" with a new coding scheme.
$ws.Range("A24").Value = "FinalTx_coll"
$ws.Range("B24").Value = "0 = BPAP,  1 = ASV,  2 = CPAP,  3 = Other,  "

# "PercOSA" (row 25) keeps its name but gets a new coding scheme.
$ws.Range("A25").Value = "PercOSA"
$ws.Range("B25").Value = "0 = mostly_OSA,  1 = mostly_CSA,  "

# "StudyType" (shifted up from row 27 to row 26) is unchanged in content.
$ws.Range("A26").Value = "StudyType"
$ws.Range("B26").Value = "0 = hst,  1 = psg,  "
